$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45243
$ws.Range("M2").Value = 50
$ws.Range("R2").Value = "Región Metropolitana"

# Row 3
$ws.Range("D3").Value = 44179
$ws.Range("M3").Value = 40
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("R3").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S3").Value = 3000

# Row 4
$ws.Range("D4").Value = 45244
$ws.Range("M4").Value = 70

# Row 5
$ws.Range("D5").Value = 45239
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 35000
$ws.Range("O5").Value = 35000
$ws.Range("P5").Value = 35000
$ws.Range("S5").Value = 7000
